$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 79.20169066666666
$ws.Cells.Item(2, 8).Value = 237.605072
$ws.Cells.Item(2, 9).Value = 0.1882902756436698
$ws.Cells.Item(2, 10).Value = 0.1882902756436699
$ws.Cells.Item(2, 13).Value = 6.045145666666667
$ws.Cells.Item(2, 14).Value = 18.135437
$ws.Cells.Item(2, 15).Value = 0.8160840232643366
$ws.Cells.Item(2, 16).Value = 0.8160840232643367
$ws.Cells.Item(2, 17).Value = 478.7857571262738
$ws.Cells.Item(2, 18).Value = 4309.071814136465
$ws.Cells.Item(2, 19).Value = 0.153660685688837
$ws.Cells.Item(2, 20).Value = 0.1536606856888371

$ws.Cells.Item(3, 7).Value = 79.20169066666666
$ws.Cells.Item(3, 8).Value = 237.605072
$ws.Cells.Item(3, 9).Value = 0.1882902756436698
$ws.Cells.Item(3, 10).Value = 0.1882902756436699
$ws.Cells.Item(3, 15).Value = 0.09212864864242169
$ws.Cells.Item(3, 16).Value = 0.09212864864242169
$ws.Cells.Item(3, 17).Value = 54.05066578419556
$ws.Cells.Item(3, 18).Value = 486.4559920577601
$ws.Cells.Item(3, 19).Value = 0.01734692864756039
$ws.Cells.Item(3, 20).Value = 0.01734692864756039

$ws.Cells.Item(4, 7).Value = 79.20169066666666
$ws.Cells.Item(4, 8).Value = 237.605072
$ws.Cells.Item(4, 9).Value = 0.1882902756436698
$ws.Cells.Item(4, 10).Value = 0.1882902756436699
$ws.Cells.Item(4, 13).Value = 0.6799149999999999
$ws.Cells.Item(4, 14).Value = 2.039745
$ws.Cells.Item(4, 15).Value = 0.09178732809324164
$ws.Cells.Item(4, 16).Value = 0.09178732809324165
$ws.Cells.Item(4, 17).Value = 53.85041750962666
$ws.Cells.Item(4, 18).Value = 484.65375758664
$ws.Cells.Item(4, 19).Value = 0.01728266130727243
$ws.Cells.Item(4, 20).Value = 0.01728266130727243

$ws.Cells.Item(5, 9).Value = 0.3031383606299999
$ws.Cells.Item(5, 10).Value = 0.3031383606299999
$ws.Cells.Item(5, 13).Value = 6.045145666666667
$ws.Cells.Item(5, 14).Value = 18.135437
$ws.Cells.Item(5, 15).Value = 0.8160840232643366
$ws.Cells.Item(5, 16).Value = 0.8160840232643367
$ws.Cells.Item(5, 17).Value = 770.8222265440789
$ws.Cells.Item(5, 18).Value = 6937.400038896711
$ws.Cells.Item(5, 19).Value = 0.2473863729486857
$ws.Cells.Item(5, 20).Value = 0.2473863729486858

$ws.Cells.Item(6, 9).Value = 0.3031383606299999
$ws.Cells.Item(6, 10).Value = 0.3031383606299999
$ws.Cells.Item(6, 15).Value = 0.09212864864242169
$ws.Cells.Item(6, 16).Value = 0.09212864864242169
$ws.Cells.Item(6, 19).Value = 0.02792772751652098
$ws.Cells.Item(6, 20).Value = 0.02792772751652098

$ws.Cells.Item(7, 9).Value = 0.3031383606299999
$ws.Cells.Item(7, 10).Value = 0.3031383606299999
$ws.Cells.Item(7, 13).Value = 0.6799149999999999
$ws.Cells.Item(7, 14).Value = 2.039745
$ws.Cells.Item(7, 15).Value = 0.09178732809324164
$ws.Cells.Item(7, 16).Value = 0.09178732809324165
$ws.Cells.Item(7, 17).Value = 86.69660303648332
$ws.Cells.Item(7, 18).Value = 780.2694273283499
$ws.Cells.Item(7, 19).Value = 0.02782426016479321
$ws.Cells.Item(7, 20).Value = 0.02782426016479321

$ws.Cells.Item(8, 7).Value = 128.6091306666667
$ws.Cells.Item(8, 8).Value = 385.827392
$ws.Cells.Item(8, 9).Value = 0.3057491381773125
$ws.Cells.Item(8, 10).Value = 0.3057491381773124
$ws.Cells.Item(8, 13).Value = 6.045145666666667
$ws.Cells.Item(8, 14).Value = 18.135437
$ws.Cells.Item(8, 15).Value = 0.8160840232643366
$ws.Cells.Item(8, 16).Value = 0.8160840232643367
$ws.Cells.Item(8, 17).Value = 777.4609289433673
$ws.Cells.Item(8, 18).Value = 6997.148360490306
$ws.Cells.Item(8, 19).Value = 0.2495169867933447
$ws.Cells.Item(8, 20).Value = 0.2495169867933447

$ws.Cells.Item(9, 7).Value = 128.6091306666667
$ws.Cells.Item(9, 8).Value = 385.827392
$ws.Cells.Item(9, 9).Value = 0.3057491381773125
$ws.Cells.Item(9, 10).Value = 0.3057491381773124
$ws.Cells.Item(9, 15).Value = 0.09212864864242169
$ws.Cells.Item(9, 16).Value = 0.09212864864242169
$ws.Cells.Item(9, 17).Value = 87.76844382926224
$ws.Cells.Item(9, 18).Value = 789.9159944633601
$ws.Cells.Item(9, 19).Value = 0.02816825492386086
$ws.Cells.Item(9, 20).Value = 0.02816825492386086

$ws.Cells.Item(10, 7).Value = 128.6091306666667
$ws.Cells.Item(10, 8).Value = 385.827392
$ws.Cells.Item(10, 9).Value = 0.3057491381773125
$ws.Cells.Item(10, 10).Value = 0.3057491381773124
$ws.Cells.Item(10, 13).Value = 0.6799149999999999
$ws.Cells.Item(10, 14).Value = 2.039745
$ws.Cells.Item(10, 15).Value = 0.09178732809324164
$ws.Cells.Item(10, 16).Value = 0.09178732809324165
$ws.Cells.Item(10, 17).Value = 87.44327707722667
$ws.Cells.Item(10, 18).Value = 786.98949369504
$ws.Cells.Item(10, 19).Value = 0.02806389646010685
$ws.Cells.Item(10, 20).Value = 0.02806389646010685

$ws.Cells.Item(11, 7).Value = 85.31435366666666
$ws.Cells.Item(11, 8).Value = 255.943061
$ws.Cells.Item(11, 9).Value = 0.2028222255490178
$ws.Cells.Item(11, 10).Value = 0.2028222255490178
$ws.Cells.Item(11, 13).Value = 6.045145666666667
$ws.Cells.Item(11, 14).Value = 18.135437
$ws.Cells.Item(11, 15).Value = 0.8160840232643366
$ws.Cells.Item(11, 16).Value = 0.8160840232643367
$ws.Cells.Item(11, 17).Value = 515.7376953725175
$ws.Cells.Item(11, 18).Value = 4641.639258352658
$ws.Cells.Item(11, 19).Value = 0.1655199778334691
$ws.Cells.Item(11, 20).Value = 0.1655199778334692

$ws.Cells.Item(12, 7).Value = 85.31435366666666
$ws.Cells.Item(12, 8).Value = 255.943061
$ws.Cells.Item(12, 9).Value = 0.2028222255490178
$ws.Cells.Item(12, 10).Value = 0.2028222255490178
$ws.Cells.Item(12, 15).Value = 0.09212864864242169
$ws.Cells.Item(12, 16).Value = 0.09212864864242169
$ws.Cells.Item(12, 17).Value = 58.22221189745889
$ws.Cells.Item(12, 18).Value = 523.99990707713
$ws.Cells.Item(12, 19).Value = 0.01868573755447946
$ws.Cells.Item(12, 20).Value = 0.01868573755447946

$ws.Cells.Item(13, 7).Value = 85.31435366666666
$ws.Cells.Item(13, 8).Value = 255.943061
$ws.Cells.Item(13, 9).Value = 0.2028222255490178
$ws.Cells.Item(13, 10).Value = 0.2028222255490178
$ws.Cells.Item(13, 13).Value = 0.6799149999999999
$ws.Cells.Item(13, 14).Value = 2.039745
$ws.Cells.Item(13, 15).Value = 0.09178732809324164
$ws.Cells.Item(13, 16).Value = 0.09178732809324165
$ws.Cells.Item(13, 17).Value = 58.00650877327165
$ws.Cells.Item(13, 18).Value = 522.058578959445
$ws.Cells.Item(13, 19).Value = 0.01861651016106915
$ws.Cells.Item(13, 20).Value = 0.01861651016106915
